$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BL (index 64) is a new "data pull" column added after BK (index 63),
# mirroring the existing per-timestamp price columns C..BK.

# 1) Header cell BL1: new timestamp label, styled like the other header cells (BK1).
$headerSrc = $ws.Cells.Item(1, 63)
$headerDst = $ws.Cells.Item(1, 64)
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)
$headerDst.Value = "2025-07-10T10:45"

# 2) Data rows: numeric values for column BL, row by row.
$blValues = @{
    2 = 100
    3 = 400
    4 = 500
    5 = 600
    6 = 400
    7 = 15000
    8 = 1000
    9 = 1200
    10 = 1600
    11 = 8300
    12 = 10400
    13 = 59900
    14 = 875000
    15 = 250000
    16 = 59999000
    17 = 65500
    18 = 350000
    19 = 32000
    20 = 4975000
    21 = 167400
    22 = 30000000
    23 = 400000
    24 = 148800
    25 = 2398000
    26 = 295000
    27 = 3750000
    28 = 549800
    29 = 175000
    30 = 500000
    31 = 1500000
    32 = 648000
    33 = 10000
    34 = 43900
    35 = 100000
    36 = 70000
    37 = 65000000
    38 = 1000000000
    39 = 1400000
    40 = 1349999800
    41 = 39900
    42 = 235000000
    43 = 940000000
    44 = 38999000
    45 = 24999900
    46 = 320000000
    47 = 278500
    48 = 200
    49 = 200
    50 = 200
    51 = 400
    52 = 1500
    53 = 600
    54 = 1800
    55 = 20000
    56 = 7200
    57 = 125000
    58 = 11900
    59 = 58000
    60 = 51000
    61 = 245000
    62 = 800000
    63 = 2799900
    64 = 1450000
    65 = 3950000
    66 = 39999900
    67 = 900000
    68 = 100
    69 = 200
    70 = 500
    71 = 1300
    72 = 600
    73 = 1500
    74 = 1900
    75 = 1100
    76 = 10600
    77 = 98500
    78 = 924900
    79 = 22999900
    80 = 125000000
    81 = 15000000000
    82 = 740000000
    83 = 3499500
    84 = 1400
    85 = 400000
    86 = 7750000000
    87 = 65000000000
    88 = 400000000000
    89 = 998000000
    90 = 1000000000
    91 = 20000000000
    92 = 1500000000
    94 = 300
    95 = 200
    96 = 300
    97 = 500
    98 = 500
    99 = 600
    100 = 1000
    101 = 134500
    102 = 50000
    103 = 526100
    104 = 599900
    105 = 10800
    106 = 3900
    107 = 11900
    108 = 49900
    109 = 2500000
    110 = 548999000
    111 = 20000000
    112 = 2000000000
    113 = 18000000000
    114 = 95000000000
    115 = 59000000000
    116 = 20000
    117 = 285000000
    118 = 330000000
    119 = 100000000
    121 = 12000
    122 = 538900
    123 = 200000
    124 = 1100000
    125 = 800000
    126 = 8000000
    127 = 40000000
    129 = 59300
    130 = 54700
    131 = 99000000
    132 = 29999900
    133 = 6000000
    134 = 69999000
    135 = 80000000
    136 = 390000000000
    137 = 45000000
    138 = 19000000
    139 = 3000000000
    140 = 12000000000
    141 = 24900000000
    142 = 7999999900
    143 = 245000000
    144 = 75000000000
    145 = 55000000000
    146 = 60000000000
    147 = 65000000000
    148 = 39500000000
    149 = 10300000000
    150 = 18500000000
    151 = 44000000000
    152 = 85000000000
    153 = 3000000
    154 = 1995000000
    155 = 5200000000
    156 = 1000000000
    157 = 838900000000
    158 = 2500000000
    159 = 13999999000
    160 = 315000000000
    161 = 990000000000
    162 = 1000000000000
    164 = 999900000000
    165 = 50000000000
    166 = 15000000000
}

foreach ($r in $blValues.Keys) {
    $ws.Cells.Item($r, 64).Value = $blValues[$r]
}

# 3) Rows that are blank placeholders across the whole data block (matching the
#    existing blank BK cells in those rows) get a blank BL cell too, formatted
#    like the neighboring blank cell so the row stays rectangular.
$blankRows = @(93, 120, 128, 163)
foreach ($r in $blankRows) {
    $src = $ws.Cells.Item($r, 63)
    $dst = $ws.Cells.Item($r, 64)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}
